# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# Rename the "Requested quantity" header to more descriptive, metric-specific
# names on the existing sheets.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add a new worksheet with the PO forecast (ds / yhat / yhat_lower / yhat_upper
# style output) as the third, trailing tab.
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

# Header row - match the bold/bordered/centered header styling used on the
# other two sheets.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"
$wsForecast.Range("A1:D1").Font.Bold = $true
$wsForecast.Range("A1:D1").HorizontalAlignment = -4108
$wsForecast.Range("A1:D1").VerticalAlignment = -4160
$wsForecast.Range("A1:D1").Borders.LineStyle = 1

# Data rows: ds (date serial), PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(45158.99999999999, 46, -39.07081884572379, 127.4443950536063),
    @(45193.99999999999, 62, -20.28361089380866, 146.1919916878238),
    @(45214.99999999999, 71, -12.66869360942976, 160.2563542430465),
    @(45221.99999999999, 75, -12.27626362383687, 158.3710679556844),
    @(45228.99999999999, 78, -2.290807041919751, 167.4641038538775),
    @(45235.99999999999, 81, -1.082325752679096, 162.8047072176871),
    @(45242.99999999999, 84, 0.8484076718976091, 166.1782898770095),
    @(45249.99999999999, 87, 2.209025691664855, 179.1413093129543),
    @(45277.99999999999, 100, 15.77075566999813, 185.9685324313849),
    @(45298.99999999999, 110, 24.35130387574558, 196.7298259068668),
    @(45305.99999999999, 113, 31.47596603627312, 190.2511533564908),
    @(45312.99999999999, 116, 27.61957836492062, 202.7846194573669),
    @(45326.99999999999, 123, 37.39069393502497, 203.9082197871821),
    @(45333.99999999999, 126, 41.91665298396944, 205.1682424965823),
    @(45340.99999999999, 129, 51.10059771625389, 212.3848700030416),
    @(45347.99999999999, 132, 47.37596506673712, 216.6359751392795),
    @(45354.99999999999, 135, 53.52541695182018, 218.7787982713472),
    @(45361.99999999999, 139, 60.80526038111317, 222.7180490171536),
    @(45368.99999999999, 142, 50.2181217786487, 225.5680713372261),
    @(45375.99999999999, 145, 59.41127375943513, 232.420114132411),
    @(45382.99999999999, 148, 59.1078219349513, 230.2962388318446)
)

$row = 2
foreach ($item in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $item[0]
    $wsForecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($row, 2).Value = $item[1]
    $wsForecast.Cells.Item($row, 3).Value = $item[2]
    $wsForecast.Cells.Item($row, 4).Value = $item[3]
    $row++
}

# Put the new sheet after "Monthly Trend" (i.e. as the last tab), matching
# the order it was inserted into the workbook.
$lastSheet = $wb.Worksheets.Item("Monthly Trend")
$wsForecast.Move($null, $lastSheet)
